$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc04cb9d38eaee60b1e31c4040511b41c5707c8e/e2e/9434839e-203f-4168-9621-5775ee0c037a.md"
$displayName = "9434839e-203f-4168-9621-5775ee0c037a.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4194360096d42eaa3ce46c68aaa5b01b7f30cf60/e2e/9434839e-203f-4168-9621-5775ee0c037a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc04cb9d38eaee60b1e31c4040511b41c5707c8e/e2e/9434839e-203f-4168-9621-5775ee0c037a.md."

# --- zh-cn sheet : row 7 (9434839e-... file) ---
# Latest Target File (J7) becomes the same handoff xlf file already shown in G7
$zhcn.Range("J7").Value = $zhcn.Range("G7").Text
# Latest Handback DateTime (K7)
$zhcn.Range("K7").Value = "2016-08-17 10:52:27"
# Error Detail (P7)
$zhcn.Range("P7").Value = $errorMessage
# Latest Handback File (I7) gets a hyperlink pointing at the handback markdown file
$zhcn.Hyperlinks.Add($zhcn.Range("I7"), $targetUrl, "", "", $displayName)

# --- de-de sheet : row 7 (9434839e-... file) ---
$dede.Range("J7").Value = $dede.Range("G7").Text
$dede.Range("K7").Value = "2016-08-17 10:52:35"
$dede.Range("P7").Value = $errorMessage
$dede.Hyperlinks.Add($dede.Range("I7"), $targetUrl, "", "", $displayName)
